# Update template: move the "Comment" column (J) to the end, shifting
# "SamplePortion" (K) and "SamplePortionUnit" (L) earlier.
# Net effect: new J = old K, new K = old L, new L = old J, for rows 1-5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($cell, $text)

    if ($text -eq "" -or $text -eq $null) {
        # Writing a plain "" via Value2 makes Excel blank out (remove) the
        # cell entirely. Force a real, stored empty-string cell instead by
        # entering it as literal text (leading apostrophe) and then
        # stripping the resulting "quote prefix" formatting so the cell is
        # left clean.
        $cell.Value2 = "'"
        $cell.ClearFormats()
    } else {
        $cell.Value2 = $text
    }
}

for ($row = 1; $row -le 5; $row++) {
    $jCell = $ws.Cells.Item($row, 10)   # column J
    $kCell = $ws.Cells.Item($row, 11)   # column K
    $lCell = $ws.Cells.Item($row, 12)   # column L

    $oldJ = $jCell.Value2
    $oldK = $kCell.Value2
    $oldL = $lCell.Value2

    Set-CellText $jCell $oldK
    Set-CellText $kCell $oldL
    Set-CellText $lCell $oldJ
}
